$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 19
$ws.Range("AA5").Value = 7.5
$ws.Range("AH5").Value = 41
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 3.4
$ws.Range("G6").Value = 3.1
$ws.Range("I6").Value = 2.6
$ws.Range("AA10").Value = 6.3
$ws.Range("AB10").Value = 16.5
$ws.Range("AC10").Value = 90
$ws.Range("AD10").Value = 800
$ws.Range("AE10").Value = 6.6
$ws.Range("AF10").Value = 9.5
$ws.Range("AH10").Value = 19.5
$ws.Range("AI10").Value = 18.5
$ws.Range("AJ10").Value = 35
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 2.12
$ws.Range("M10").Value = 2.62
$ws.Range("R10").Value = 1.87
$ws.Range("S10").Value = 1.75
$ws.Range("V10").Value = 11.75
$ws.Range("AA14").Value = 6.8
$ws.Range("AE14").Value = 11
$ws.Range("AF14").Value = 20
$ws.Range("AG14").Value = 12.5
$ws.Range("AH14").Value = 55
$ws.Range("AI14").Value = 32
$ws.Range("G14").Value = 1.9
$ws.Range("I14").Value = 3.65
$ws.Range("L14").Value = 1.27
$ws.Range("M14").Value = 3.1
$ws.Range("N14").Value = 1.8
$ws.Range("O14").Value = 1.8
$ws.Range("S14").Value = 1.91
$ws.Range("T14").Value = 7.5
$ws.Range("U14").Value = 9.25
$ws.Range("W14").Value = 16.5
$ws.Range("X14").Value = 14.5
$ws.Range("L15").Value = 1.25
$ws.Range("M15").Value = 3.75
$ws.Range("N15").Value = 1.8
$ws.Range("O15").Value = 2
$ws.Range("AB16").Value = 21
$ws.Range("AE16").Value = 13.5
$ws.Range("K16").Value = 7.6
$ws.Range("L16").Value = 1.32
$ws.Range("M16").Value = 3.2
$ws.Range("N16").Value = 1.95
$ws.Range("T16").Value = 5.8
$ws.Range("U16").Value = 6.9
$ws.Range("X16").Value = 14
$ws.Range("Z16").Value = 7.6
$ws.Range("AC18").Value = 51
$ws.Range("G18").Value = 1.53
$ws.Range("H18").Value = 4.1
$ws.Range("I18").Value = 5.75
$ws.Range("N18").Value = 1.7
$ws.Range("O18").Value = 2.1
$ws.Range("R18").Value = 1.8
$ws.Range("S18").Value = 1.91
$ws.Range("W18").Value = 11
$ws.Range("AA19").Value = 10
$ws.Range("AC19").Value = 81
$ws.Range("AE19").Value = 6
$ws.Range("AG19").Value = 9
$ws.Range("AH19").Value = 8.5
$ws.Range("AI19").Value = 13
$ws.Range("G19").Value = 5.75
$ws.Range("H19").Value = 4.75
$ws.Range("I19").Value = 1.36
$ws.Range("L19").Value = 1.25
$ws.Range("M19").Value = 3.25
$ws.Range("N19").Value = 1.88
$ws.Range("O19").Value = 1.93
$ws.Range("P19").Value = 1.38
$ws.Range("Q19").Value = 2.6
$ws.Range("R19").Value = 2.07
$ws.Range("S19").Value = 1.6
$ws.Range("T19").Value = 15
$ws.Range("U19").Value = 34
$ws.Range("V19").Value = 21
$ws.Range("W19").Value = 81
$ws.Range("Z19").Value = 11
$ws.Range("AI20").Value = 26
$ws.Range("G20").Value = 2
$ws.Range("N20").Value = 1.73
$ws.Range("O20").Value = 2.08
$ws.Range("AB21").Value = 12
$ws.Range("AD21").Value = 126
$ws.Range("AE21").Value = 11
$ws.Range("AF21").Value = 15
$ws.Range("AJ21").Value = 23
$ws.Range("H21").Value = 3.6
$ws.Range("L21").Value = 1.2
$ws.Range("M21").Value = 4.33
$ws.Range("P21").Value = 1.3
$ws.Range("Q21").Value = 3.4
$ws.Range("R21").Value = 1.57
$ws.Range("S21").Value = 2.25
$ws.Range("T21").Value = 11
$ws.Range("U21").Value = 15
$ws.Range("Y21").Value = 23
$ws.Range("J22").Value = 1.02
$ws.Range("K22").Value = 19
$ws.Range("R26").Value = 1.67
$ws.Range("AE32").Value = 10.5
$ws.Range("AF32").Value = 13
$ws.Range("AH32").Value = 24
$ws.Range("AJ32").Value = 22
$ws.Range("G32").Value = 2.7
$ws.Range("I32").Value = 2.27
$ws.Range("N32").Value = 1.57
$ws.Range("O32").Value = 2.12
$ws.Range("R32").Value = 1.5
$ws.Range("S32").Value = 2.27
$ws.Range("U32").Value = 16.5
$ws.Range("X32").Value = 20
$ws.Range("AD36").Value = 301
$ws.Range("AF36").Value = 29
$ws.Range("T36").Value = 7
$ws.Range("Y36").Value = 26
$ws.Range("Z36").Value = 11
$ws.Range("K37").Value = 13
$ws.Range("P37").Value = 1.3
$ws.Range("AI38").Value = 41
$ws.Range("G38").Value = 1.5
$ws.Range("H38").Value = 3.9
$ws.Range("I38").Value = 6
$ws.Range("AA39").Value = 8.5
$ws.Range("AB39").Value = 19
$ws.Range("AC39").Value = 90
$ws.Range("AD39").Value = 700
$ws.Range("AE39").Value = 7
$ws.Range("AF39").Value = 6.7
$ws.Range("AG39").Value = 8.25
$ws.Range("AH39").Value = 9.5
$ws.Range("AI39").Value = 11.5
$ws.Range("AJ39").Value = 27
$ws.Range("G39").Value = 6.2
$ws.Range("H39").Value = 4.3
$ws.Range("I39").Value = 1.44
$ws.Range("M39").Value = 3.4
$ws.Range("R39").Value = 1.87
$ws.Range("S39").Value = 1.75
$ws.Range("T39").Value = 16.5
$ws.Range("U39").Value = 40
$ws.Range("V39").Value = 20
$ws.Range("W39").Value = 120
$ws.Range("X39").Value = 70
$ws.Range("Y39").Value = 65
$ws.Range("AD40").Value = 450
$ws.Range("AH40").Value = 45
$ws.Range("G40").Value = 1.98
$ws.Range("I40").Value = 3.3
$ws.Range("O40").Value = 1.82
$ws.Range("V40").Value = 8.5
$ws.Range("T42").Value = 9.5
$ws.Range("AB43").Value = 15
$ws.Range("AC43").Value = 65
$ws.Range("AE43").Value = 14.5
$ws.Range("AF43").Value = 29
$ws.Range("AG43").Value = 15
$ws.Range("G43").Value = 1.62
$ws.Range("I43").Value = 4.75
$ws.Range("T43").Value = 7.2
$ws.Range("U43").Value = 7.8
$ws.Range("X43").Value = 13
$ws.Range("Y43").Value = 25
$ws.Range("AA47").Value = 6.4
$ws.Range("AB47").Value = 15.5
$ws.Range("AC47").Value = 80
$ws.Range("AD47").Value = 700
$ws.Range("AE47").Value = 9.25
$ws.Range("AF47").Value = 17
$ws.Range("AG47").Value = 11.75
$ws.Range("AH47").Value = 45
$ws.Range("AI47").Value = 32
$ws.Range("AJ47").Value = 40
$ws.Range("G47").Value = 2.05
$ws.Range("H47").Value = 3.3
$ws.Range("I47").Value = 3.3
$ws.Range("L47").Value = 1.34
$ws.Range("M47").Value = 2.75
$ws.Range("N47").Value = 2
$ws.Range("O47").Value = 1.65
$ws.Range("P47").Value = 1.44
$ws.Range("Q47").Value = 2.42
$ws.Range("R47").Value = 1.82
$ws.Range("S47").Value = 1.78
$ws.Range("T47").Value = 6.7
$ws.Range("U47").Value = 9.25
$ws.Range("V47").Value = 8.75
$ws.Range("W47").Value = 18.5
$ws.Range("X47").Value = 18
$ws.Range("Y47").Value = 32
$ws.Range("Z47").Value = 8.75
$ws.Range("L48").Value = 1.29
$ws.Range("M48").Value = 3.5
$ws.Range("N48").Value = 1.93
$ws.Range("O48").Value = 1.88
$ws.Range("P48").Value = 1.33
$ws.Range("P57").Value = 1.41
$ws.Range("Q57").Value = 2.62
$ws.Range("L63").Value = 1.29
$ws.Range("M63").Value = 3.5
$ws.Range("N63").Value = 1.9
$ws.Range("O63").Value = 1.9
$ws.Range("AA64").Value = 6.4
$ws.Range("AB64").Value = 25
$ws.Range("AE64").Value = 4.75
$ws.Range("AF64").Value = 7.3
$ws.Range("AH64").Value = 16.5
$ws.Range("AI64").Value = 22
$ws.Range("G64").Value = 4.15
$ws.Range("H64").Value = 3
$ws.Range("I64").Value = 1.95
$ws.Range("M64").Value = 2.07
$ws.Range("O64").Value = 1.36
$ws.Range("Q64").Value = 2.07
$ws.Range("R64").Value = 2.37
$ws.Range("S64").Value = 1.45
$ws.Range("T64").Value = 7.8
$ws.Range("U64").Value = 20
$ws.Range("V64").Value = 16.5
$ws.Range("W64").Value = 75
$ws.Range("X64").Value = 60
$ws.Range("Z64").Value = 4.8

Write-Host "Applied 229 cell updates"
